$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 0.5979736666666666
$ws.Range("H2").Value = 1.793921
$ws.Range("I2").Value = 0.03342655292740804
$ws.Range("J2").Value = 0.03342655292740804
$ws.Range("M2").Value = 1.492477333333333
$ws.Range("N2").Value = 4.477432
$ws.Range("O2").Value = 0.02769484181536182
$ws.Range("P2").Value = 0.02769484181536182
$ws.Range("Q2").Value = 0.8924621434302222
$ws.Range("R2").Value = 8.032159290872
$ws.Range("S2").Value = 0.0009257430957573856
$ws.Range("T2").Value = 0.0009257430957573856

# Row 3
$ws.Range("G3").Value = 0.5979736666666666
$ws.Range("H3").Value = 1.793921
$ws.Range("I3").Value = 0.03342655292740804
$ws.Range("J3").Value = 0.03342655292740804
$ws.Range("O3").Value = 0.6282762845978157
$ws.Range("P3").Value = 0.6282762845978156
$ws.Range("Q3").Value = 20.24610948698489
$ws.Range("R3").Value = 182.214985382864
$ws.Range("S3").Value = 0.02100111048014417
$ws.Range("T3").Value = 0.02100111048014416

# Row 4
$ws.Range("G4").Value = 0.5979736666666666
$ws.Range("H4").Value = 1.793921
$ws.Range("I4").Value = 0.03342655292740804
$ws.Range("J4").Value = 0.03342655292740804
$ws.Range("N4").Value = 55.619234
$ws.Range("O4").Value = 0.3440288735868225
$ws.Range("P4").Value = 0.3440288735868225
$ws.Range("Q4").Value = 11.08627909739044
$ws.Range("R4").Value = 99.77651187651399
$ws.Range("S4").Value = 0.0114996993515065
$ws.Range("T4").Value = 0.0114996993515065

# Row 5
$ws.Range("I5").Value = 0.8874158839838097
$ws.Range("J5").Value = 0.8874158839838097
$ws.Range("M5").Value = 1.492477333333333
$ws.Range("N5").Value = 4.477432
$ws.Range("O5").Value = 0.02769484181536182
$ws.Range("P5").Value = 0.02769484181536182
$ws.Range("Q5").Value = 23.69329208591022
$ws.Range("R5").Value = 213.239628773192
$ws.Range("S5").Value = 0.02457684253137109
$ws.Range("T5").Value = 0.02457684253137109

# Row 6
$ws.Range("I6").Value = 0.8874158839838097
$ws.Range("J6").Value = 0.8874158839838097
$ws.Range("O6").Value = 0.6282762845978157
$ws.Range("P6").Value = 0.6282762845978156
$ws.Range("S6").Value = 0.5575423544824342
$ws.Range("T6").Value = 0.5575423544824341

# Row 7
$ws.Range("I7").Value = 0.8874158839838097
$ws.Range("J7").Value = 0.8874158839838097
$ws.Range("N7").Value = 55.619234
$ws.Range("O7").Value = 0.3440288735868225
$ws.Range("P7").Value = 0.3440288735868225
$ws.Range("S7").Value = 0.3052966869700044
$ws.Range("T7").Value = 0.3052966869700044

# Row 8
$ws.Range("I8").Value = 0.07915756308878232
$ws.Range("J8").Value = 0.07915756308878232
$ws.Range("M8").Value = 1.492477333333333
$ws.Range("N8").Value = 4.477432
$ws.Range("O8").Value = 0.02769484181536182
$ws.Range("P8").Value = 0.02769484181536182
$ws.Range("Q8").Value = 2.113443422549334
$ws.Range("R8").Value = 19.020990802944
$ws.Range("S8").Value = 0.00219225618823335
$ws.Range("T8").Value = 0.00219225618823335

# Row 9
$ws.Range("I9").Value = 0.07915756308878232
$ws.Range("J9").Value = 0.07915756308878232
$ws.Range("O9").Value = 0.6282762845978157
$ws.Range("P9").Value = 0.6282762845978156
$ws.Range("Q9").Value = 47.94489855112535
$ws.Range("R9").Value = 431.5040869601281
$ws.Range("S9").Value = 0.04973281963523735
$ws.Range("T9").Value = 0.04973281963523734

# Row 10
$ws.Range("I10").Value = 0.07915756308878232
$ws.Range("J10").Value = 0.07915756308878232
$ws.Range("N10").Value = 55.619234
$ws.Range("O10").Value = 0.3440288735868225
$ws.Range("P10").Value = 0.3440288735868225
$ws.Range("Q10").Value = 26.25346499165867
$ws.Range("S10").Value = 0.02723248726531162
$ws.Range("T10").Value = 0.02723248726531162
